# Week 13 logging update
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# YDS sheet: append this week's run/pass yardage logs (space separated lists)
# ---------------------------------------------------------------------------
$wsYDS = $wb.Worksheets.Item("YDS")

$wsYDS.Range("B2").Value = $wsYDS.Range("B2").Value2 + " 3 7 1 3 0 17 7 3 15 -1 6 2 1 1 10 1"
$wsYDS.Range("C2").Value = $wsYDS.Range("C2").Value2 + " 22 1 3 1 16 2 6 2 2 4 1 9 5 -4 -3 13 4 1 0 -6 7 2 3 8 1 6 0 4 3"
$wsYDS.Range("B3").Value = $wsYDS.Range("B3").Value2 + " 1 8 9 3 10 5 3 9 7 5 0 18 10 -5 5 2 15 9 34 7 3 28 14 13 16 2 11 7"
$wsYDS.Range("C3").Value = $wsYDS.Range("C3").Value2 + " 23 -1 13 6 7 4 7 6 5 6 5 6 0 10 35 10 3 4 10 9 10 12 6"

# ---------------------------------------------------------------------------
# OFF sheet totals
# ---------------------------------------------------------------------------
$wsOFF = $wb.Worksheets.Item("OFF")

$wsOFF.Range("C2").Value = 334
$wsOFF.Range("D2").Value = 32
$wsOFF.Range("E2").Value = 12
$wsOFF.Range("F2").Value = 80
$wsOFF.Range("G2").Value = 104
$wsOFF.Range("L2").Value = 552
$wsOFF.Range("M2").Value = 382
$wsOFF.Range("Q2").Value = 963

$wsOFF.Range("C3").Value = 316
$wsOFF.Range("E3").Value = 45
$wsOFF.Range("F3").Value = 204
$wsOFF.Range("G3").Value = 74
$wsOFF.Range("H3").Value = 32
$wsOFF.Range("I3").Value = 95
$wsOFF.Range("J3").Value = 94
$wsOFF.Range("N3").Value = 27

# ---------------------------------------------------------------------------
# DEF sheet totals
# ---------------------------------------------------------------------------
$wsDEF = $wb.Worksheets.Item("DEF")

$wsDEF.Range("C2").Value = 353
$wsDEF.Range("F2").Value = 107
$wsDEF.Range("G2").Value = 110
$wsDEF.Range("H2").Value = 7
$wsDEF.Range("I2").Value = 10
$wsDEF.Range("J2").Value = 49
$wsDEF.Range("L2").Value = 545
$wsDEF.Range("M2").Value = 368
$wsDEF.Range("Q2").Value = 978

$wsDEF.Range("B3").Value = 16
$wsDEF.Range("C3").Value = 297
$wsDEF.Range("E3").Value = 59
$wsDEF.Range("F3").Value = 183
$wsDEF.Range("H3").Value = 50
$wsDEF.Range("I3").Value = 105
$wsDEF.Range("J3").Value = 91
$wsDEF.Range("N3").Value = 20

# ---------------------------------------------------------------------------
# ST sheet (special teams) totals + append kick/return distance logs
# ---------------------------------------------------------------------------
$wsST = $wb.Worksheets.Item("ST")

$wsST.Range("B2").Value = 154
$wsST.Range("D2").Value = 90
$wsST.Range("H2").Value = 3
$wsST.Range("J2").Value = 64
$wsST.Range("K2").Value = 61
$wsST.Range("N2").Value = 15
$wsST.Range("O2").Value = 12

$wsST.Range("B3").Value = 97

$wsST.Range("D3").Value = $wsST.Range("D3").Value2 + " 42 67 38 55"
$wsST.Range("D4").Value = $wsST.Range("D4").Value2 + " 0 16 0 13"
$wsST.Range("D5").Value = $wsST.Range("D5").Value2 + " 6 0 0 0"
$wsST.Range("B6").Value = $wsST.Range("B6").Value2 + " 19 24"

# ---------------------------------------------------------------------------
# TURNS sheet
# ---------------------------------------------------------------------------
$wsTURNS = $wb.Worksheets.Item("TURNS")

$wsTURNS.Range("C2").Value = 7
$wsTURNS.Range("D3").Value = 9

# ---------------------------------------------------------------------------
# PEN sheet
# ---------------------------------------------------------------------------
$wsPEN = $wb.Worksheets.Item("PEN")

$wsPEN.Range("D2").Value = 21
$wsPEN.Range("D3").Value = 13
